# Auto-generated edit script: update F-column 'want to go' counts
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 42016
$ws.Range("F5").Value = 9644
$ws.Range("F7").Value = 898
$ws.Range("F8").Value = 921
$ws.Range("F9").Value = 737
$ws.Range("F12").Value = 304
$ws.Range("F13").Value = 924
$ws.Range("F14").Value = 87
$ws.Range("F16").Value = 749
$ws.Range("F17").Value = 321
$ws.Range("F18").Value = 1448
$ws.Range("F20").Value = 687
$ws.Range("F21").Value = 716
$ws.Range("F23").Value = 692
$ws.Range("F24").Value = 751
$ws.Range("F27").Value = 63
$ws.Range("F28").Value = 510
$ws.Range("F29").Value = 536
$ws.Range("F30").Value = 58
$ws.Range("F31").Value = 248
$ws.Range("F32").Value = 937
$ws.Range("F33").Value = 20
$ws.Range("F34").Value = 438
$ws.Range("F35").Value = 103
$ws.Range("F37").Value = 148
$ws.Range("F38").Value = 410
$ws.Range("F39").Value = 1295
$ws.Range("F40").Value = 302
$ws.Range("F41").Value = 1260
$ws.Range("F42").Value = 379
$ws.Range("F43").Value = 98
$ws.Range("F45").Value = 39
$ws.Range("F46").Value = 35
$ws.Range("F49").Value = 65

$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 40
$ws.Range("F3").Value = 216
$ws.Range("F5").Value = 4454
$ws.Range("F7").Value = 338
$ws.Range("F10").Value = 79
$ws.Range("F11").Value = 132
$ws.Range("F15").Value = 23
$ws.Range("F17").Value = 165
$ws.Range("F19").Value = 4386

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 2042
$ws.Range("F3").Value = 533
$ws.Range("F4").Value = 423

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2042
$ws.Range("F3").Value = 533
$ws.Range("F5").Value = 40
$ws.Range("F6").Value = 216
$ws.Range("F7").Value = 338
$ws.Range("F9").Value = 9644
$ws.Range("F11").Value = 898
$ws.Range("F12").Value = 898
$ws.Range("F13").Value = 79
$ws.Range("F14").Value = 423
$ws.Range("F15").Value = 921
$ws.Range("F16").Value = 132
$ws.Range("F18").Value = 304
$ws.Range("F19").Value = 924
$ws.Range("F21").Value = 87
$ws.Range("F23").Value = 749
$ws.Range("F24").Value = 321
$ws.Range("F25").Value = 1448
$ws.Range("F27").Value = 687
$ws.Range("F28").Value = 716
$ws.Range("F30").Value = 692
$ws.Range("F31").Value = 751
$ws.Range("F33").Value = 63
$ws.Range("F34").Value = 510
$ws.Range("F35").Value = 58
$ws.Range("F36").Value = 248
$ws.Range("F37").Value = 937
$ws.Range("F39").Value = 20
$ws.Range("F40").Value = 438
$ws.Range("F41").Value = 103
$ws.Range("F43").Value = 1260
$ws.Range("F44").Value = 379
$ws.Range("F45").Value = 98
$ws.Range("F46").Value = 39
$ws.Range("F50").Value = 65
